$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add missing value to row 36 (Pruebas Positivas) ---
$ws.Range("C36").Value = 27

# --- Append new daily rows 37-44 ---
# Column A uses the existing date-style formatting (m/d/yy -> built-in date format)
$dateFmt = "m/d/yy"

$ws.Range("A37").NumberFormat = $dateFmt
$ws.Range("A37").Value = 43956

$ws.Range("A38").NumberFormat = $dateFmt
$ws.Range("A38").Value = 43957

$ws.Range("A39").NumberFormat = $dateFmt
$ws.Range("A39").Value = 43958

$ws.Range("A40").NumberFormat = $dateFmt
$ws.Range("A40").Value = 43959
$ws.Range("B40").Value = 815
$ws.Range("C40").Value = 68

$ws.Range("A41").NumberFormat = $dateFmt
$ws.Range("A41").Value = 43960
$ws.Range("B41").Value = 967
$ws.Range("C41").Value = 67

$ws.Range("A42").NumberFormat = $dateFmt
$ws.Range("A42").Value = 43961

$ws.Range("A43").NumberFormat = $dateFmt
$ws.Range("A43").Value = 43962

$ws.Range("A44").NumberFormat = $dateFmt
$ws.Range("A44").Value = 43963

# --- Grow the table to cover the new rows ---
$table = $ws.ListObjects.Item("Condicion_Pacientes")
$table.Resize($ws.Range("A1:F44"))

# --- Update sheet view to match the scrolled/selected state after the edit ---
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("A45").Select()
